$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D11").Value = "[1, 0, 0, 0, 1, 0, 0]"
$ws.Range("E11").Value = "['Normal', 'RegulationViolation']"

$ws.Range("D25").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E25").Value = "['Normal']"

$ws.Range("D26").Value = "[0, 0, 1, 0, 0, 0, 1]"
$ws.Range("E26").Value = "['HardwareFault', 'SoftwareFault']"

$ws.Range("D35").Value = "[1, 0, 1, 0, 0, 0, 0]"
$ws.Range("E35").Value = "['Normal', 'HardwareFault']"

$ws.Range("D38").Value = "[0, 0, 1, 0, 0, 0, 0]"
$ws.Range("E38").Value = "['HardwareFault']"

$ws.Range("D39").Value = "[1, 0, 1, 0, 0, 0, 0]"
$ws.Range("E39").Value = "['Normal', 'HardwareFault']"

$ws.Range("D69").Value = "[1, 1, 0, 0, 0, 0, 0]"
$ws.Range("E69").Value = "['Normal', 'SurroundingEnvironment']"

$ws.Range("D73").Value = "[1, 0, 0, 1, 0, 0, 0]"
$ws.Range("E73").Value = "['Normal', 'ParamViolation']"

$ws.Range("D74").Value = "[1, 0, 0, 0, 0, 0, 0]"
$ws.Range("E74").Value = "['Normal']"

$ws.Range("D81").Value = "[1, 0, 1, 0, 0, 0, 0]"
$ws.Range("E81").Value = "['Normal', 'HardwareFault']"

$ws.Range("D92").Value = "[1, 0, 1, 0, 0, 0, 1]"
$ws.Range("E92").Value = "['Normal', 'HardwareFault', 'SoftwareFault']"
